$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Pass/Fail" result column (F) for the test steps in rows 2-6.
# Row 2 and 4 passed outright ("Pass"); rows 3, 5 and 6 are the lower-case "pass".
$ws.Range("F2").Value = "Pass"
$ws.Range("F3").Value = "pass"
$ws.Range("F4").Value = "Pass"
$ws.Range("F5").Value = "pass"
$ws.Range("F6").Value = "pass"

# Match the author's final view state: scrolled one column right, with F6 selected.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F6").Select()
